# Weekly refresh: a new price observation is inserted as row 32 (with a new
# date), and all the previously-recorded rows 32-99 shift down by one row to
# rows 33-100. Row 32's non-date columns reuse the data that used to sit in
# (old) row 33, i.e. effectively: new A32:T32 = old A33:T33 except for the
# date (column D), which becomes the new observation's date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current block that needs to move (old rows 32-99, all columns
# A:T) before writing anything, then shift the whole block down by one row.
$block = $ws.Range("A32:T99").Value2
$ws.Range("A33:T100").Value2 = $block

# The brand new row 100 doesn't inherit the date-formatted style that column
# D carries throughout the table, so re-apply it explicitly.
$ws.Range("D100").NumberFormat = $ws.Range("D99").NumberFormat

# Row 32 becomes the new observation: same data as the old row 33 (which has
# now also been copied down into row 34) but tagged with the new date.
$ws.Range("A32:T32").Value2 = $ws.Range("A34:T34").Value2
$ws.Range("D32").Value2 = 44614
